$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "VAR"
$ws.Range("A1").Select() | Out-Null
